$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = '$2b$10$w7UTEHxgiHzfqskR92FNou8bJeESigklW6kK9JIc4eNtyOMlOQfJe'
$ws.Range("C3").Value = '$2b$10$Aa2MpvDQyP6I/Dw.Qlmls.gRrBEvaYIjQDImGd5SVAb4CUHiClpHa'
$ws.Range("C4").Value = '$2b$10$woMEqSoLET5oX1ICEIBsYeJavfuNTr6b9U1ZHPeKNJVdEj.GcyQvO'
$ws.Range("C5").Value = '$2b$10$.JPeyM5ikTU0Dj.jw82WAurxydU4ERxoL2bPnxZJMUq8QpToB6w8i'
$ws.Range("C6").Value = '$2b$10$ztsEgCrp/pVTKIrDt0y.wOfAMTaw6Z9ybaojpwtf/jVOuUr0ghY96'
$ws.Range("C7").Value = '$2b$10$5P1J.CRB7jnHdiYOxEJkkOJFSwOKR7MgGPUojFuZa3jmfrq2whKgm'
$ws.Range("C8").Value = '$2b$10$EOcFMmNzsLFhia1Ma08rJeUR4g5hCTikNxVUyAR1rQrwi9QVCfquK'
$ws.Range("C9").Value = '$2b$10$8Kpw7EG1GfCbV/6EVIfSzuXWn5sd4PjqwIFpQsUUgRlaSj9C66.NC'
$ws.Range("C10").Value = '$2b$10$YFiIaIX051JhnhEM7I9NgOttLbXDm.Hr3.QUC20N.F/DDrDeeZH2W'
$ws.Range("C11").Value = '$2b$10$uwavtLAhcGFKnzfEh16clukOeMlPYq6JhkQ5KRMrstqWu.pvTLBbq'
$ws.Range("C12").Value = '$2b$10$8KkmWvXtQzRUZ1V1kMDkxerJdXPPnojKvKvjPu1CocpEVumRBCR32'
$ws.Range("C13").Value = '$2b$10$kKDLNSsAXqVAbmiwcH/0lOvJl.e52xH2pnrxGRRaQGi3fIo27ALK2'
$ws.Range("C14").Value = '$2b$10$Y4/wy.h4i3PaFDGfV9ljZeUVnquhQwxNYn1e3NfmW6ozAwGfh0yfS'
$ws.Range("C15").Value = '$2b$10$BU7oDlB6mkXjaqubbwMUlunso2SvuWvGIN9RHBJPPwFsLa4g9NagW'
$ws.Range("C16").Value = '$2b$10$DAyWU/3qTySd53/2e7A64.3zA1b4F7QCQ0Z7ZfFT/DMWYmwfBlSnm'
$ws.Range("C17").Value = '$2b$10$uXh5ccczGw7Q5zwDO92xBucPb72l6S18KEhVcZlMQlvs902XF6gbC'
$ws.Range("C18").Value = '$2b$10$ezo2285gmT69KM9P9SxER.PcLYay.796cv0k2S7j4vu.e1yJjhgbK'
$ws.Range("C19").Value = '$2b$10$yy/z0AGAbovY1Lt9cjyOT.Uhavmkh8GS1j37J3bmn0qesWXm9qnlS'
$ws.Range("C20").Value = '$2b$10$iNgQTYuEwpxG4PfmnUMXtudcKyp1gtIkdnoXx/ZPHiM3lfQunPOpi'
$ws.Range("C21").Value = '$2b$10$yyinjloyJgBYl9Fo2t17QO4ghZRFvALC8JeDg0ksujUpzCvN37PRC'
$ws.Range("C22").Value = '$2b$10$F8E5beGdLbfDOxLDiwe/.efNXH9de/ez74qp3yxCvasFWupB/jo4a'
$ws.Range("C23").Value = '$2b$10$iVpEnCqAg9jWeEH40keEm.4ejSuvOKfBBYZM4jXr6zU5YRYgqw.gS'
$ws.Range("C24").Value = '$2b$10$Z8JHTNBZbvvpqjOM/hrekubfP8zF32qKU2o6pR9MG0pK7QakdwCWq'
$ws.Range("C25").Value = '$2b$10$yXd0/wsXMHTXOXIQiC63U.dVd2WIbsE.IXQBEuzKHAT/hbbhrUVe6'
$ws.Range("C26").Value = '$2b$10$z9XUl/AcEfTH9zZVGPS2led2PwcVFrmnUp1SosRs0Vhf7vRenmjIe'
$ws.Range("C27").Value = '$2b$10$b4zdlfaPQXb/uAPaLSH/aObglLLEEkxfTRrhhHpwpaLwHSC5tMWnW'
$ws.Range("C28").Value = '$2b$10$RBN0PSBGbYHAGfsL4.LISenI7dY1iL1weIXVhY1rZXU91GLrpU.ca'
$ws.Range("C29").Value = '$2b$10$nmQOHH0MeCunWxFrE6EiheQP/e4PLzUSYOz60s34N8HHqmR.keani'
$ws.Range("C30").Value = '$2b$10$bK1Gnhzi2L6/hiF6GpLZE.e.VjhdNtr3pp9dVJmMRlSGJFibp5NNW'
$ws.Range("C31").Value = '$2b$10$o5pPROJmnQ.DZbg.LVvtKe8EZ6dXDxDtl8hvkqkFDZUCNafJX3LGW'